{"js": "// Load all paragraphs in the body so we can locate the title, the old\n// \"Play Chupacabra...\" promo paragraph and the promo/meta paragraph that\n// follows it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// 1) Insert a new \"Meta description\" paragraph directly after the H1 title\n//    paragraph (paragraph 0): a bold \"Meta description\" label run followed\n//    by a normal run with the description text (prefixed with \": \").\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\n\nconst labelRange = metaPara.insertText(\"Meta description\", Word.InsertLocation.end);\nlabelRange.font.bold = true;\n\nconst descriptionText =\n  \": Explore the mystery of Chupacabra with this high-volatility online slot game! \" +\n  \"Play for free now and earn exciting multipliers and free spins.\";\nconst descRange = metaPara.insertText(descriptionText, Word.InsertLocation.end);\ndescRange.font.bold = false;\n\nawait context.sync();\n\n// 2) Near the end of the document, remove the bold paragraph that repeats\n//    the title (\"Play Chupacabra Free Now - A Thrilling Online Slot Game\")\n//    and replace the text of the following italic paragraph (previously the\n//    meta description, now an AI image-generation prompt) while keeping its\n//    italic formatting intact.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items,text\");\nawait context.sync();\n\nlet promoParagraph = null;\nlet imagePromptParagraph = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const para = paragraphs2.items[i];\n  const text = para.text.trim();\n  if (i > 0 && text === \"Play Chupacabra Free Now - A Thrilling Online Slot Game\") {\n    promoParagraph = para;\n  } else if (text.indexOf(\"Explore the mystery of Chupacabra\") === 0) {\n    imagePromptParagraph = para;\n  }\n}\n\nif (promoParagraph) {\n  promoParagraph.delete();\n  await context.sync();\n}\n\nif (imagePromptParagraph) {\n  const newImagePromptText =\n    \"Create a feature image for Chupacabra: Design a fun and exciting cartoon-style image for the Chupacabra online slot game. \" +\n    \"The focal point of the image should be a happy Maya warrior, who is donning glasses. \" +\n    \"The warrior should be depicted in a dynamic pose, with their arms outstretched and a big smile on their face. \" +\n    \"The background of the image should feature a mysterious jungle scene, complete with vines, trees, and exotic plants. \" +\n    \"The image should be bright and full of color, with a mix of greens and blues to create a sense of depth and excitement. \" +\n    \"The overall effect should be engaging and fun, with the Maya warrior looking like they are ready for an epic adventure in the world of Chupacabra.\";\n  imagePromptParagraph.insertText(newImagePromptText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph directly after the H1\n#    title paragraph (paragraph 1): a bold \"Meta description\" label\n#    followed by a normal run with the description text (prefixed with\n#    \": \").\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titleRng = $titlePara.Range.Duplicate\n$titleRng.Collapse(0)              # wdCollapseEnd\n$titleRng.InsertParagraphAfter() | Out-Null\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n\n$metaLabel = \"Meta description\"\n$metaDescription = \": Explore the mystery of Chupacabra with this high-volatility online slot game! Play for free now and earn exciting multipliers and free spins.\"\n\n$metaRng = $metaPara.Range\n$metaRng.InsertBefore($metaLabel + $metaDescription) | Out-Null\n\n# Bold only the \"Meta description\" label, leaving the rest of the text\n# (the description itself) in normal weight.\n$metaPara = $d.Paragraphs.Item(2)\n$labelRng = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $metaLabel.Length)\n$labelRng.Bold = 1\n\n# ---------------------------------------------------------------------\n# 2) Near the end of the document, remove the bold paragraph that\n#    repeats the title (\"Play Chupacabra Free Now - A Thrilling Online\n#    Slot Game\") and replace the text of the following italic paragraph\n#    (previously the meta description, now an AI image-generation\n#    prompt) while keeping its italic formatting intact.\n# ---------------------------------------------------------------------\n$promoParaIndex = -1\n$imagePromptParaIndex = -1\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx++\n    $text = $p.Range.Text.Trim()\n    if ($idx -ne 1 -and $text -eq \"Play Chupacabra Free Now - A Thrilling Online Slot Game\") {\n        $promoParaIndex = $idx\n    } elseif ($text.StartsWith(\"Explore the mystery of Chupacabra\")) {\n        $imagePromptParaIndex = $idx\n    }\n}\n\nif ($promoParaIndex -gt 0) {\n    $promoPara = $d.Paragraphs.Item($promoParaIndex)\n    $promoPara.Range.Delete() | Out-Null\n    if ($imagePromptParaIndex -gt $promoParaIndex) {\n        $imagePromptParaIndex = $imagePromptParaIndex - 1\n    }\n}\n\nif ($imagePromptParaIndex -gt 0) {\n    $imagePromptPara = $d.Paragraphs.Item($imagePromptParaIndex)\n    $imgRng = $imagePromptPara.Range.Duplicate\n    $imgRng.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark\n    $imgRng.Text = \"Create a feature image for Chupacabra: Design a fun and exciting cartoon-style image for the Chupacabra online slot game. The focal point of the image should be a happy Maya warrior, who is donning glasses. The warrior should be depicted in a dynamic pose, with their arms outstretched and a big smile on their face. The background of the image should feature a mysterious jungle scene, complete with vines, trees, and exotic plants. The image should be bright and full of color, with a mix of greens and blues to create a sense of depth and excitement. The overall effect should be engaging and fun, with the Maya warrior looking like they are ready for an epic adventure in the world of Chupacabra.\"\n}\n"}
